$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 91
$ws.Range("I9").Value = 95.55556
$ws.Range("K9").Value = 95.55556
$ws.Range("M9").Value = 73.44444
$ws.Range("H88").Value = 2017.0646
$ws.Range("I88").Value = 1094.5
$ws.Range("J88").Value = 2337.9565
$ws.Range("K88").Value = 1094.5
$ws.Range("L88").Value = 2337.9565
$ws.Range("M88").Value = -688.5
$ws.Range("N88").Value = -3149.9565
$ws.Range("H91").Value = 2017.0646
$ws.Range("I91").Value = 1094.5
$ws.Range("J91").Value = 2337.9565
$ws.Range("K91").Value = 1094.5
$ws.Range("L91").Value = 2337.9565
$ws.Range("M91").Value = 309.5
$ws.Range("N91").Value = -5145.9565
$ws.Range("H129").Value = 1160.174
$ws.Range("J129").Value = 1256.2683
$ws.Range("L129").Value = 3768.8049
$ws.Range("N129").Value = -13768.8049
$ws.Range("H137").Value = 3243.3572
$ws.Range("I137").Value = 3630.158
$ws.Range("J137").Value = 2426.7778
$ws.Range("K137").Value = 10890.474
$ws.Range("L137").Value = 7280.3334
$ws.Range("M137").Value = -8340.474
$ws.Range("N137").Value = -12380.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1730511
$ws.Range("I2").Value = 402.875
$ws.Range("K2").Value = 402.875
$ws.Range("M2").Value = -289.875
$ws.Range("H61").Value = 2340.195
$ws.Range("I61").Value = 2055.1667
$ws.Range("J61").Value = 3117.5454
$ws.Range("K61").Value = 2055.1667
$ws.Range("L61").Value = 3117.5454
$ws.Range("M61").Value = -1843.1667
$ws.Range("N61").Value = -3541.5454
$ws.Range("H74").Value = 2328717.8
$ws.Range("I74").Value = 3125505.5
$ws.Range("J74").Value = 10790
$ws.Range("K74").Value = 3125505.5
$ws.Range("L74").Value = 10790
$ws.Range("M74").Value = -3124631.5
$ws.Range("N74").Value = -12538
$ws.Range("H77").Value = 2328717.8
$ws.Range("I77").Value = 3125505.5
$ws.Range("J77").Value = 10790
$ws.Range("K77").Value = 15627527.5
$ws.Range("L77").Value = 53950
$ws.Range("M77").Value = -15623159.5
$ws.Range("N77").Value = -62686
$ws.Range("H116").Value = 1730511
$ws.Range("I116").Value = 402.875
$ws.Range("K116").Value = 402.875
$ws.Range("M116").Value = 1891.125
$ws.Range("H136").Value = 2340.195
$ws.Range("I136").Value = 2055.1667
$ws.Range("J136").Value = 3117.5454
$ws.Range("K136").Value = 6165.500100000001
$ws.Range("L136").Value = 9352.636200000001
$ws.Range("M136").Value = -3615.500100000001
$ws.Range("N136").Value = -14452.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1730511
$ws.Range("I3").Value = 402.875
$ws.Range("K3").Value = 402.875
$ws.Range("M3").Value = -288.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 463.1842
$ws.Range("I22").Value = 472.25
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 472.25
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -122.25
$ws.Range("N22").Value = -1000
$ws.Range("H31").Value = 4066886.2
$ws.Range("I31").Value = 1453.9259
$ws.Range("J31").Value = 11907363
$ws.Range("K31").Value = 1453.9259
$ws.Range("L31").Value = 11907363
$ws.Range("M31").Value = -1158.9259
$ws.Range("N31").Value = -11907953
$ws.Range("H34").Value = 4066886.2
$ws.Range("I34").Value = 1453.9259
$ws.Range("J34").Value = 11907363
$ws.Range("K34").Value = 1453.9259
$ws.Range("L34").Value = 11907363
$ws.Range("M34").Value = -1251.9259
$ws.Range("N34").Value = -11907767
$ws.Range("H86").Value = 3526.7896
$ws.Range("I86").Value = 3045.2727
$ws.Range("J86").Value = 4188.875
$ws.Range("K86").Value = 3045.2727
$ws.Range("L86").Value = 4188.875
$ws.Range("M86").Value = -1922.2727
$ws.Range("N86").Value = -6434.875
$ws.Range("H89").Value = 3526.7896
$ws.Range("I89").Value = 3045.2727
$ws.Range("J89").Value = 4188.875
$ws.Range("K89").Value = 15226.3635
$ws.Range("L89").Value = 20944.375
$ws.Range("M89").Value = -9610.363499999999
$ws.Range("N89").Value = -32176.375
$ws.Range("H99").Value = 1399.75
$ws.Range("I99").Value = 1366.5
$ws.Range("J99").Value = 1499.5
$ws.Range("K99").Value = 1366.5
$ws.Range("L99").Value = 1499.5
$ws.Range("M99").Value = 131.5
$ws.Range("N99").Value = -4495.5
$ws.Range("H126").Value = 1399.75
$ws.Range("I126").Value = 1366.5
$ws.Range("J126").Value = 1499.5
$ws.Range("K126").Value = 4099.5
$ws.Range("L126").Value = 4498.5
$ws.Range("M126").Value = -1629.5
$ws.Range("N126").Value = -9438.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2176384.5
$ws.Range("J131").Value = 3336626.5
$ws.Range("L131").Value = 10009879.5
$ws.Range("N131").Value = -10019959.5
$ws.Range("H132").Value = 52632636
$ws.Range("I132").Value = 90910344
$ws.Range("J132").Value = 795.5
$ws.Range("K132").Value = 818193096
$ws.Range("L132").Value = 7159.5
$ws.Range("M132").Value = -818190566
$ws.Range("N132").Value = -12219.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1779.174
$ws.Range("I16").Value = 1922.9048
$ws.Range("J16").Value = 270
$ws.Range("K16").Value = 1922.9048
$ws.Range("L16").Value = 270
$ws.Range("M16").Value = -1752.9048
$ws.Range("N16").Value = -610
$ws.Range("H22").Value = 409.83334
$ws.Range("I22").Value = 312.25
$ws.Range("J22").Value = 605
$ws.Range("K22").Value = 312.25
$ws.Range("L22").Value = 605
$ws.Range("M22").Value = -17.25
$ws.Range("N22").Value = -1195
$ws.Range("H27").Value = 409.83334
$ws.Range("I27").Value = 312.25
$ws.Range("J27").Value = 605
$ws.Range("K27").Value = 312.25
$ws.Range("L27").Value = 605
$ws.Range("M27").Value = -205.25
$ws.Range("N27").Value = -819
$ws.Range("H31").Value = 1502.3334
$ws.Range("I31").Value = 1004.6667
$ws.Range("J31").Value = 2000
$ws.Range("K31").Value = 1004.6667
$ws.Range("L31").Value = 2000
$ws.Range("M31").Value = -756.6667
$ws.Range("N31").Value = -2496
$ws.Range("H82").Value = 2858.1
$ws.Range("I82").Value = 2837.375
$ws.Range("J82").Value = 2871.9167
$ws.Range("K82").Value = 2837.375
$ws.Range("L82").Value = 2871.9167
$ws.Range("M82").Value = -2476.375
$ws.Range("N82").Value = -3593.9167
$ws.Range("H85").Value = 2858.1
$ws.Range("I85").Value = 2837.375
$ws.Range("J85").Value = 2871.9167
$ws.Range("K85").Value = 2837.375
$ws.Range("L85").Value = 2871.9167
$ws.Range("M85").Value = -1589.375
$ws.Range("N85").Value = -5367.9167
$ws.Range("H100").Value = 2107.3572
$ws.Range("I100").Value = 1980.3
$ws.Range("J100").Value = 2425
$ws.Range("K100").Value = 1980.3
$ws.Range("L100").Value = 2425
$ws.Range("M100").Value = -1439.3
$ws.Range("N100").Value = -3507
$ws.Range("H122").Value = 12235.818
$ws.Range("I122").Value = 17543.428
$ws.Range("J122").Value = 2947.5
$ws.Range("K122").Value = 52630.284
$ws.Range("L122").Value = 8842.5
$ws.Range("M122").Value = -50180.284
$ws.Range("N122").Value = -13742.5
$ws.Range("H136").Value = 2226.6206
$ws.Range("I136").Value = 1580.6818
$ws.Range("J136").Value = 4256.7144
$ws.Range("K136").Value = 4742.0454
$ws.Range("L136").Value = 12770.1432
$ws.Range("M136").Value = -2192.0454
$ws.Range("N136").Value = -17870.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 11000
$ws.Range("I61").Value = 5000
$ws.Range("J61").Value = 12500
$ws.Range("K61").Value = 5000
$ws.Range("L61").Value = 12500
$ws.Range("M61").Value = -4708
$ws.Range("N61").Value = -13084
$ws.Range("H81").Value = 3399.25
$ws.Range("I81").Value = 1600
$ws.Range("J81").Value = 4298.875
$ws.Range("K81").Value = 3200
$ws.Range("L81").Value = 8597.75
$ws.Range("M81").Value = -2139
$ws.Range("N81").Value = -10719.75
$ws.Range("H84").Value = 3399.25
$ws.Range("I84").Value = 1600
$ws.Range("J84").Value = 4298.875
$ws.Range("K84").Value = 16000
$ws.Range("L84").Value = 42988.75
$ws.Range("M84").Value = -10696
$ws.Range("N84").Value = -53596.75
$ws.Range("H107").Value = 307.875
$ws.Range("I107").Value = 261.64285
$ws.Range("J107").Value = 372.6
$ws.Range("K107").Value = 784.9285500000001
$ws.Range("L107").Value = 1117.8
$ws.Range("M107").Value = 1135.07145
$ws.Range("N107").Value = -4957.8
$ws.Range("H113").Value = 581.5
$ws.Range("I113").Value = 459.8
$ws.Range("K113").Value = 1379.4
$ws.Range("M113").Value = 790.5999999999999
$ws.Range("H132").Value = 6131411.5
$ws.Range("I132").Value = 7693849
$ws.Range("J132").Value = 1849
$ws.Range("K132").Value = 23081547
$ws.Range("L132").Value = 5547
$ws.Range("M132").Value = -23079017
$ws.Range("N132").Value = -10607
